# "found mistakes were fixed"
#
# The lintroller_test sheet lists boolean flags for lintroller_81 .. lintroller_99
# in rows 2-20, with row 21 left blank. The edit:
#   1. Adds a new entry "lintroller_100" with bool value 1.
#   2. Because "lintroller_100" sorts before "lintroller_81" .. "lintroller_99"
#      (plain text/alphabetic ordering), the new row is inserted at the TOP of the
#      data (row 2), pushing every existing row down by one - which is what fills
#      the previously-empty row 21 with what used to be row 20 (lintroller_99).
#      Every existing name keeps its original bool value; nothing else changes.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row right above the current first data row, shifting
# lintroller_81..lintroller_99 (and the blank row 21) down by one.
$ws.Rows(2).Insert()

# Populate the newly inserted row with the new entry.
$ws.Range("A2").Value = "lintroller_100"
$ws.Range("B2").Value = 1

# The insert duplicated the last (blank) row as a new trailing row 22;
# drop it so the used range goes back to A1:B21.
$ws.Rows(22).Delete()

# The freshly written cell comes in unformatted (default style) rather than
# inheriting the sheet's normal data style - normalize B2's font so it keeps
# that default look.
$ws.Range("B2").Font.Bold = $false

# Reflect the author's final cursor position in the saved view.
[void]$ws.Range("D5").Select()
